$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.036.40"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.30%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.906.30"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.85%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9994"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "333.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.76%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9997"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.10%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4637"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.37%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4078"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.79"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.76%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08001"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.20%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.005"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.41%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.73"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.57%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.911.81"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.44%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.945"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.73%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.111"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.73%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "89.09"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.34%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9982"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.39%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001035"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.74%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06540"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.13%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.49"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.31%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.002"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.12%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "29.013.27"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.20%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.459"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.11%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.16"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.03%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.243"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.10%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.133.49"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.16%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "157.96"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.41%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.73"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.17%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.106"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.12%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.390"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.24%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "119.14"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.75%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9831"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.41%  "

$ws.Range("E33").Value = "  -1.08%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.419"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.05%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.599"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.74%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.303"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.85%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06086"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.07%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02230"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.05%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.397"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.27%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.166"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.67%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5820"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.09%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9991"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.12%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.18"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.51%  "

$ws.Range("E44").Value = "  -2.43%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.243"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.93%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.338"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +14.91%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.16"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.03%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5506"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.46%  "

$ws.Range("E49").Value = "  -2.65%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07032"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.45%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "47.29"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +22.53%  "
